$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "K" (strikeouts) column values, replacing old "Strike#" totals
$ws.Range("G2").Value = 4
$ws.Range("G3").Value = 4
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 4
$ws.Range("G7").Value = 2
$ws.Range("G8").Value = 2
$ws.Range("G9").Value = 9
$ws.Range("G10").Value = 5
$ws.Range("G11").Value = 1
$ws.Range("G13").Value = 3
$ws.Range("G14").Value = 0
